$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the time-interval labels in column C (rows 7-12) with the new
# continuing sequence of 5-minute slots, swapping rows 8 and 9 relative to
# the simple sequential order (matches the shared-string index swap in the
# source diff).
$ws.Range("C7").Value = "21:50-21:55"
$ws.Range("C9").Value = "22:0-22:5"
$ws.Range("C8").Value = "21:55-22:0"
$ws.Range("C10").Value = "22:5-22:10"
$ws.Range("C11").Value = "22:10-22:15"
$ws.Range("C12").Value = "22:15-22:20"

# Move the active selection from C16 to C15.
$ws.Range("C15").Select()
